$wb = $excel.ActiveWorkbook

# Sheets that contain this data table: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F9").Value = 6421
    $ws.Range("F12").Value = 315
    $ws.Range("F13").Value = 104
    $ws.Range("F16").Value = 6410
}

$wb.Save()
